# Add files via upload
# - Rename header N1 from "syst_u" to "syst_tot"
# - Add seven new per-source systematic-uncertainty columns (Q:W) with
#   headers syst0_c .. syst6_c and their numeric values for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells Q1:W1 ------------------------------------------------
$headers = @("syst0_c", "syst1_c", "syst2_c", "syst3_c", "syst4_c", "syst5_c", "syst6_c")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, 17 + $i)
    $cell.Value = $headers[$i]
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# --- New data cells Q2:W11 --------------------------------------------------
$data = @{
  2  = @(0.00007000000000000001, 0.00018, 0.00001, 0.00012, 0.00006, 0.00107, 0.0005)
  3  = @(0.00005, 0.00036, 0.00006, 0.00008, 0.00005, 0.00129, 0.0013)
  4  = @(0.00029, 0.00046, 0.00007000000000000001, 0.00013, 0.00004, 0.00151, 0.0006)
  5  = @(0.00049, 0.00065, 0.00012, 0.00039, 0.00003, 0.00165, 0.0011)
  6  = @(0.00047, 0.00089, 0.00012, 0.00046, 0.00004, 0.00134, 0.0009)
  7  = @(0.0005099999999999999, 0.00078, 0.00014, 0.00053, 0.00002, 0.00251, 0.0022)
  8  = @(0.00057, 0.00058, 0.00006, 0.00042, 0.00002, 0.00187, 0.0017)
  9  = @(0.00055, 0.00048, 0.00013, 0.00073, 0.00005, 0.00183, 0.0017)
  10 = @(0.0003, 0.00005, 0.00047, 0.00082, 0.00031, 0.00312, 0.002)
  11 = @(0.00037, 0.00008999999999999999, 0.00048, 0.00149, 0.00049, 0.00385, 0.0004)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 17 + $i).Value = $vals[$i]
    }
}

# --- Column N header: syst_u -> syst_tot -----------------------------------
$ws.Cells.Item(1, 14).Value = "syst_tot"
